$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mohiuddin")

$ws.Range("C13").Value = "Completed the soldering of all three boards. `nTested a simple sender receiver on them. `nStudied some more products and included them in comparison sheet.`nhttps://github.com/frenziopen/FrenziTech/blob/main/Documentation/Comparison%20Sheet.xlsx`nAttached images of the products in comparison sheet."

$ws.Range("C16").Value = "Studied ESP32 pinout in detail from https://docs.espressif.com/projects/esp-idf/en/latest/esp32/hw-reference/esp32/get-started-devkitc.html#get-started-esp32-devkitc-board-front `nand created pins utilization diagram uploaded FrenziTech/README.md at main ` frenziopen/FrenziTech (github.com)."

$ws.Range("C13").WrapText = $true
$ws.Range("C16").WrapText = $true

$ws.Range("C17").Select()
